$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated data values from the Preprocessing and Cleaning script re-run ---
$ws.Range("F9").Value = "Caribbean"
$ws.Range("G25").Value = 7001762
$ws.Range("O25").Value = 3.979313911590529
$ws.Range("U25").Value = 0.1428211927226318
$ws.Range("F32").Value = "Unknown"
$ws.Range("G47").Value = 106766
$ws.Range("O47").Value = 593.1444444444444
$ws.Range("U47").Value = 18.73255530786955
$ws.Range("F51").Value = "Unknown"
$ws.Range("H53").Value = 133.505341239138
$ws.Range("O53").Value = 15.51997830774862
$ws.Range("J69").Value = "Unknown"
$ws.Range("J81").Value = '["Mamoudzou"]'
$ws.Range("J90").Value = "Unknown"
$ws.Range("G124").Value = 367830
$ws.Range("O124").Value = 3.571165048543689
$ws.Range("U124").Value = 2.718647201152706
$ws.Range("H134").Value = 513120
$ws.Range("O134").Value = 136.030515279077
$ws.Range("H169").Value = 323639.1380288619
$ws.Range("O169").Value = 300.7627062438889
$ws.Range("G213").Value = 9537642
$ws.Range("O213").Value = 66.65018867924529
$ws.Range("U213").Value = 0.2096954362514341
$ws.Range("J214").Value = '["Valletta"]'
$ws.Range("H241").Value = 390757
$ws.Range("O241").Value = 38.03623991380832
$ws.Range("F242").Value = "Caribbean"

# --- Refresh the 'timestamp' column (N) for every data row with the new run's timestamps ---
$nTimestamps = @{
    2 = "2025-04-06T02:01:40.565222"
    3 = "2025-04-06T02:01:40.565222"
    4 = "2025-04-06T02:01:40.565222"
    5 = "2025-04-06T02:01:40.565222"
    6 = "2025-04-06T02:01:40.581230"
    7 = "2025-04-06T02:01:40.581230"
    8 = "2025-04-06T02:01:40.581230"
    9 = "2025-04-06T02:01:40.581230"
    10 = "2025-04-06T02:01:40.596907"
    11 = "2025-04-06T02:01:40.596907"
    12 = "2025-04-06T02:01:40.596907"
    13 = "2025-04-06T02:01:40.596907"
    14 = "2025-04-06T02:01:40.596907"
    15 = "2025-04-06T02:01:40.612525"
    16 = "2025-04-06T02:01:40.612525"
    17 = "2025-04-06T02:01:40.612525"
    18 = "2025-04-06T02:01:40.612525"
    19 = "2025-04-06T02:01:40.612525"
    20 = "2025-04-06T02:01:40.612525"
    21 = "2025-04-06T02:01:40.628157"
    22 = "2025-04-06T02:01:40.628157"
    23 = "2025-04-06T02:01:40.628157"
    24 = "2025-04-06T02:01:40.628157"
    25 = "2025-04-06T02:01:40.628157"
    26 = "2025-04-06T02:01:40.628157"
    27 = "2025-04-06T02:01:40.643775"
    28 = "2025-04-06T02:01:40.643775"
    29 = "2025-04-06T02:01:40.643775"
    30 = "2025-04-06T02:01:40.643775"
    31 = "2025-04-06T02:01:40.643775"
    32 = "2025-04-06T02:01:40.643775"
    33 = "2025-04-06T02:01:40.659408"
    34 = "2025-04-06T02:01:40.659408"
    35 = "2025-04-06T02:01:40.659408"
    36 = "2025-04-06T02:01:40.659408"
    37 = "2025-04-06T02:01:40.659408"
    38 = "2025-04-06T02:01:40.659408"
    39 = "2025-04-06T02:01:40.675026"
    40 = "2025-04-06T02:01:40.675026"
    41 = "2025-04-06T02:01:40.675026"
    42 = "2025-04-06T02:01:40.675026"
    43 = "2025-04-06T02:01:40.690652"
    44 = "2025-04-06T02:01:40.690652"
    45 = "2025-04-06T02:01:40.690652"
    46 = "2025-04-06T02:01:40.690652"
    47 = "2025-04-06T02:01:40.706279"
    48 = "2025-04-06T02:01:40.706279"
    49 = "2025-04-06T02:01:40.706279"
    50 = "2025-04-06T02:01:40.706279"
    51 = "2025-04-06T02:01:40.721905"
    52 = "2025-04-06T02:01:40.721905"
    53 = "2025-04-06T02:01:40.721905"
    54 = "2025-04-06T02:01:40.721905"
    55 = "2025-04-06T02:01:40.721905"
    56 = "2025-04-06T02:01:40.721905"
    57 = "2025-04-06T02:01:40.737527"
    58 = "2025-04-06T02:01:40.737527"
    59 = "2025-04-06T02:01:40.737527"
    60 = "2025-04-06T02:01:40.737527"
    61 = "2025-04-06T02:01:40.753151"
    62 = "2025-04-06T02:01:40.753151"
    63 = "2025-04-06T02:01:40.753151"
    64 = "2025-04-06T02:01:40.753151"
    65 = "2025-04-06T02:01:40.753151"
    66 = "2025-04-06T02:01:40.753151"
    67 = "2025-04-06T02:01:40.768779"
    68 = "2025-04-06T02:01:40.768779"
    69 = "2025-04-06T02:01:40.768779"
    70 = "2025-04-06T02:01:40.768779"
    71 = "2025-04-06T02:01:40.784402"
    72 = "2025-04-06T02:01:40.784402"
    73 = "2025-04-06T02:01:40.784402"
    74 = "2025-04-06T02:01:40.784402"
    75 = "2025-04-06T02:01:40.784402"
    76 = "2025-04-06T02:01:40.800026"
    77 = "2025-04-06T02:01:40.800026"
    78 = "2025-04-06T02:01:40.800026"
    79 = "2025-04-06T02:01:40.800026"
    80 = "2025-04-06T02:01:40.815654"
    81 = "2025-04-06T02:01:40.815654"
    82 = "2025-04-06T02:01:40.815654"
    83 = "2025-04-06T02:01:40.815654"
    84 = "2025-04-06T02:01:40.815654"
    85 = "2025-04-06T02:01:40.831294"
    86 = "2025-04-06T02:01:40.831294"
    87 = "2025-04-06T02:01:40.831294"
    88 = "2025-04-06T02:01:40.831294"
    89 = "2025-04-06T02:01:40.846901"
    90 = "2025-04-06T02:01:40.846901"
    91 = "2025-04-06T02:01:40.846901"
    92 = "2025-04-06T02:01:40.846901"
    93 = "2025-04-06T02:01:40.846901"
    94 = "2025-04-06T02:01:40.862524"
    95 = "2025-04-06T02:01:40.862524"
    96 = "2025-04-06T02:01:40.862524"
    97 = "2025-04-06T02:01:40.862524"
    98 = "2025-04-06T02:01:40.862524"
    99 = "2025-04-06T02:01:40.862524"
    100 = "2025-04-06T02:01:40.878156"
    101 = "2025-04-06T02:01:40.878156"
    102 = "2025-04-06T02:01:40.878156"
    103 = "2025-04-06T02:01:40.878156"
    104 = "2025-04-06T02:01:40.878156"
    105 = "2025-04-06T02:01:40.893791"
    106 = "2025-04-06T02:01:40.893791"
    107 = "2025-04-06T02:01:40.893791"
    108 = "2025-04-06T02:01:40.893791"
    109 = "2025-04-06T02:01:40.893791"
    110 = "2025-04-06T02:01:40.909401"
    111 = "2025-04-06T02:01:40.909401"
    112 = "2025-04-06T02:01:40.909401"
    113 = "2025-04-06T02:01:40.909401"
    114 = "2025-04-06T02:01:40.909401"
    115 = "2025-04-06T02:01:40.925130"
    116 = "2025-04-06T02:01:40.925130"
    117 = "2025-04-06T02:01:40.925130"
    118 = "2025-04-06T02:01:40.925130"
    119 = "2025-04-06T02:01:40.925130"
    120 = "2025-04-06T02:01:40.925130"
    121 = "2025-04-06T02:01:40.940656"
    122 = "2025-04-06T02:01:40.940656"
    123 = "2025-04-06T02:01:40.940656"
    124 = "2025-04-06T02:01:40.940656"
    125 = "2025-04-06T02:01:40.940656"
    126 = "2025-04-06T02:01:40.956275"
    127 = "2025-04-06T02:01:40.956275"
    128 = "2025-04-06T02:01:40.956275"
    129 = "2025-04-06T02:01:40.956275"
    130 = "2025-04-06T02:01:40.956275"
    131 = "2025-04-06T02:01:40.971904"
    132 = "2025-04-06T02:01:40.971904"
    133 = "2025-04-06T02:01:40.971904"
    134 = "2025-04-06T02:01:40.971904"
    135 = "2025-04-06T02:01:40.971904"
    136 = "2025-04-06T02:01:40.987533"
    137 = "2025-04-06T02:01:40.987533"
    138 = "2025-04-06T02:01:40.987533"
    139 = "2025-04-06T02:01:40.987533"
    140 = "2025-04-06T02:01:40.987533"
    141 = "2025-04-06T02:01:41.003152"
    142 = "2025-04-06T02:01:41.003152"
    143 = "2025-04-06T02:01:41.003152"
    144 = "2025-04-06T02:01:41.003152"
    145 = "2025-04-06T02:01:41.003152"
    146 = "2025-04-06T02:01:41.018777"
    147 = "2025-04-06T02:01:41.018777"
    148 = "2025-04-06T02:01:41.018777"
    149 = "2025-04-06T02:01:41.018777"
    150 = "2025-04-06T02:01:41.018777"
    151 = "2025-04-06T02:01:41.018777"
    152 = "2025-04-06T02:01:41.034406"
    153 = "2025-04-06T02:01:41.034406"
    154 = "2025-04-06T02:01:41.034406"
    155 = "2025-04-06T02:01:41.034406"
    156 = "2025-04-06T02:01:41.050036"
    157 = "2025-04-06T02:01:41.050036"
    158 = "2025-04-06T02:01:41.050036"
    159 = "2025-04-06T02:01:41.050036"
    160 = "2025-04-06T02:01:41.050036"
    161 = "2025-04-06T02:01:41.050036"
    162 = "2025-04-06T02:01:41.065652"
    163 = "2025-04-06T02:01:41.065652"
    164 = "2025-04-06T02:01:41.065652"
    165 = "2025-04-06T02:01:41.065652"
    166 = "2025-04-06T02:01:41.065652"
    167 = "2025-04-06T02:01:41.082017"
    168 = "2025-04-06T02:01:41.082017"
    169 = "2025-04-06T02:01:41.082017"
    170 = "2025-04-06T02:01:41.082017"
    171 = "2025-04-06T02:01:41.082017"
    172 = "2025-04-06T02:01:41.097711"
    173 = "2025-04-06T02:01:41.097711"
    174 = "2025-04-06T02:01:41.097711"
    175 = "2025-04-06T02:01:41.097711"
    176 = "2025-04-06T02:01:41.113383"
    177 = "2025-04-06T02:01:41.113383"
    178 = "2025-04-06T02:01:41.113383"
    179 = "2025-04-06T02:01:41.113383"
    180 = "2025-04-06T02:01:41.113383"
    181 = "2025-04-06T02:01:41.113383"
    182 = "2025-04-06T02:01:41.129537"
    183 = "2025-04-06T02:01:41.129537"
    184 = "2025-04-06T02:01:41.129537"
    185 = "2025-04-06T02:01:41.129537"
    186 = "2025-04-06T02:01:41.129537"
    187 = "2025-04-06T02:01:41.129537"
    188 = "2025-04-06T02:01:41.145228"
    189 = "2025-04-06T02:01:41.145228"
    190 = "2025-04-06T02:01:41.145228"
    191 = "2025-04-06T02:01:41.145228"
    192 = "2025-04-06T02:01:41.145228"
    193 = "2025-04-06T02:01:41.161285"
    194 = "2025-04-06T02:01:41.161285"
    195 = "2025-04-06T02:01:41.161285"
    196 = "2025-04-06T02:01:41.161285"
    197 = "2025-04-06T02:01:41.161285"
    198 = "2025-04-06T02:01:41.161285"
    199 = "2025-04-06T02:01:41.176965"
    200 = "2025-04-06T02:01:41.176965"
    201 = "2025-04-06T02:01:41.176965"
    202 = "2025-04-06T02:01:41.176965"
    203 = "2025-04-06T02:01:41.176965"
    204 = "2025-04-06T02:01:41.192582"
    205 = "2025-04-06T02:01:41.192582"
    206 = "2025-04-06T02:01:41.192582"
    207 = "2025-04-06T02:01:41.192582"
    208 = "2025-04-06T02:01:41.192582"
    209 = "2025-04-06T02:01:41.208207"
    210 = "2025-04-06T02:01:41.208207"
    211 = "2025-04-06T02:01:41.208207"
    212 = "2025-04-06T02:01:41.208207"
    213 = "2025-04-06T02:01:41.208207"
    214 = "2025-04-06T02:01:41.208207"
    215 = "2025-04-06T02:01:41.223832"
    216 = "2025-04-06T02:01:41.223832"
    217 = "2025-04-06T02:01:41.223832"
    218 = "2025-04-06T02:01:41.223832"
    219 = "2025-04-06T02:01:41.223832"
    220 = "2025-04-06T02:01:41.223832"
    221 = "2025-04-06T02:01:41.239460"
    222 = "2025-04-06T02:01:41.239460"
    223 = "2025-04-06T02:01:41.239460"
    224 = "2025-04-06T02:01:41.239460"
    225 = "2025-04-06T02:01:41.239460"
    226 = "2025-04-06T02:01:41.255082"
    227 = "2025-04-06T02:01:41.255082"
    228 = "2025-04-06T02:01:41.255082"
    229 = "2025-04-06T02:01:41.255082"
    230 = "2025-04-06T02:01:41.270706"
    231 = "2025-04-06T02:01:41.270706"
    232 = "2025-04-06T02:01:41.270706"
    233 = "2025-04-06T02:01:41.270706"
    234 = "2025-04-06T02:01:41.270706"
    235 = "2025-04-06T02:01:41.270706"
    236 = "2025-04-06T02:01:41.286332"
    237 = "2025-04-06T02:01:41.286332"
    238 = "2025-04-06T02:01:41.286332"
    239 = "2025-04-06T02:01:41.286332"
    240 = "2025-04-06T02:01:41.286332"
    241 = "2025-04-06T02:01:41.301959"
    242 = "2025-04-06T02:01:41.301959"
    243 = "2025-04-06T02:01:41.301959"
    244 = "2025-04-06T02:01:41.301959"
    245 = "2025-04-06T02:01:41.301959"
    246 = "2025-04-06T02:01:41.317583"
    247 = "2025-04-06T02:01:41.317583"
    248 = "2025-04-06T02:01:41.317583"
    249 = "2025-04-06T02:01:41.317583"
    250 = "2025-04-06T02:01:41.317583"
    251 = "2025-04-06T02:01:41.317583"
}

foreach ($row in $nTimestamps.Keys) {
    $ws.Cells.Item($row, 14).Value = $nTimestamps[$row]
}
